$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "VALOR MORA" total, worker count and period count summary cells
$ws.Range("E11").Value2 = 165181
$ws.Range("C13").Value2 = 2
$ws.Range("F13").Value2 = 2

# Remove the old detail rows for the workers/periods that are no longer part of
# this statement (rows 17-28), keeping the last detail row (old row 29, which
# becomes row 17) and the trailing signature rows (old rows 34-35, which
# become rows 22-23).
$ws.Rows("17:28").Delete()

# Update the remaining detail row with the new worker/period data
$ws.Range("C17").Value2 = "1143364843"
$ws.Range("D17").Value2 = "SIGRID MILEIDIS COHEN MUNEVAR"
$ws.Range("E17").Value2 = "2508"
$ws.Range("F17").Value2 = 151410
$ws.Range("G17").Value2 = 3785250

# Column D's best-fit width shrinks now that the longest name was removed
$ws.Columns("D").AutoFit()
